$wb = $excel.ActiveWorkbook
$dateFmt = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------------
# Sheet 1: "Weekly Quantity" -> append two new weekly rows (44, 45)
# ---------------------------------------------------------------------------
$wsWeekly = $wb.Worksheets.Item("Weekly Quantity")

$weeklyRows = @(44, 45)
$weeklyDates = @(45669.99999999999, 45676.99999999999)
$weeklyVals  = @(117, 12)

for ($i = 0; $i -lt $weeklyRows.Count; $i++) {
    $r = $weeklyRows[$i]
    $wsWeekly.Cells.Item($r, 1).Value = $weeklyDates[$i]
    $wsWeekly.Cells.Item($r, 1).NumberFormat = $dateFmt
    $wsWeekly.Cells.Item($r, 2).Value = $weeklyVals[$i]
}

# ---------------------------------------------------------------------------
# Sheet 2: "Monthly Trend" -> append one new monthly row (17)
# ---------------------------------------------------------------------------
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

$monthlyRows  = @(17)
$monthlyDates = @(45688.99999999999)
$monthlyVals  = @(129)

for ($i = 0; $i -lt $monthlyRows.Count; $i++) {
    $r = $monthlyRows[$i]
    $wsMonthly.Cells.Item($r, 1).Value = $monthlyDates[$i]
    $wsMonthly.Cells.Item($r, 1).NumberFormat = $dateFmt
    $wsMonthly.Cells.Item($r, 2).Value = $monthlyVals[$i]
}

# ---------------------------------------------------------------------------
# Sheet 3: "PO Forecast" -> new forecast model, rows 2-53
# (rows 2-43 keep their date, new quantity; rows 44-53 are new/shifted dates)
# ---------------------------------------------------------------------------
$wsForecast = $wb.Worksheets.Item("PO Forecast")

$forecastRows = @(
    2, 3, 4, 5, 6, 7, 8, 9, 10,
    11, 12, 13, 14, 15, 16, 17, 18, 19, 20,
    21, 22, 23, 24, 25, 26, 27, 28, 29, 30,
    31, 32, 33, 34, 35, 36, 37, 38, 39, 40,
    41, 42, 43, 44, 45, 46, 47, 48, 49, 50,
    51, 52, 53
)

$forecastDates = @(
    45137.99999999999, 45151.99999999999, 45158.99999999999, 45165.99999999999, 45172.99999999999,
    45179.99999999999, 45186.99999999999, 45193.99999999999, 45200.99999999999,
    45207.99999999999, 45214.99999999999, 45221.99999999999, 45228.99999999999, 45235.99999999999,
    45242.99999999999, 45249.99999999999, 45256.99999999999, 45263.99999999999, 45277.99999999999,
    45298.99999999999, 45312.99999999999, 45319.99999999999, 45368.99999999999, 45375.99999999999,
    45403.99999999999, 45410.99999999999, 45417.99999999999, 45424.99999999999, 45431.99999999999,
    45438.99999999999, 45445.99999999999, 45452.99999999999, 45459.99999999999, 45473.99999999999,
    45487.99999999999, 45494.99999999999, 45585.99999999999, 45592.99999999999, 45599.99999999999,
    45606.99999999999, 45634.99999999999, 45641.99999999999, 45669.99999999999, 45676.99999999999,
    45683.99999999999, 45690.99999999999, 45697.99999999999, 45704.99999999999, 45711.99999999999,
    45718.99999999999, 45725.99999999999, 45732.99999999999
)

$forecastVals = @(
    208, 209, 209, 210, 210,
    211, 211, 212, 212,
    213, 214, 214, 215, 215,
    216, 216, 217, 217, 218,
    220, 221, 222, 225, 226,
    228, 229, 229, 230, 230,
    231, 231, 232, 233, 234,
    235, 235, 242, 243, 243,
    244, 246, 247, 249, 249,
    250, 251, 251, 252, 252,
    253, 253, 254
)

for ($i = 0; $i -lt $forecastRows.Count; $i++) {
    $r = $forecastRows[$i]
    $wsForecast.Cells.Item($r, 1).Value = $forecastDates[$i]
    $wsForecast.Cells.Item($r, 1).NumberFormat = $dateFmt
    $wsForecast.Cells.Item($r, 2).Value = $forecastVals[$i]
}

Write-Host "Done"
